$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 24.50527833333333
$ws.Range("H2").Value = 73.515835
$ws.Range("I2").Value = 0.04846830138877924
$ws.Range("J2").Value = 0.04846830138877924
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.383699
$ws.Range("N2").Value = 4.151097
$ws.Range("O2").Value = 0.08080976933214185
$ws.Range("P2").Value = 0.08080976933214185
$ws.Range("Q2").Value = 33.907929124555
$ws.Range("R2").Value = 305.171362120995
$ws.Range("S2").Value = 0.00391671225514798
$ws.Range("T2").Value = 0.00391671225514798

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 24.50527833333333
$ws.Range("H3").Value = 73.515835
$ws.Range("I3").Value = 0.04846830138877924
$ws.Range("J3").Value = 0.04846830138877924
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.462094
$ws.Range("N3").Value = 7.386282
$ws.Range("O3").Value = 0.1437893994387872
$ws.Range("P3").Value = 0.1437893994387872
$ws.Range("Q3").Value = 60.33429875283
$ws.Range("R3").Value = 543.00868877547
$ws.Range("S3").Value = 0.006969227948510703
$ws.Range("T3").Value = 0.006969227948510704

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 24.50527833333333
$ws.Range("H4").Value = 73.515835
$ws.Range("I4").Value = 0.04846830138877924
$ws.Range("J4").Value = 0.04846830138877924
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 4.336036333333333
$ws.Range("N4").Value = 13.008109
$ws.Range("O4").Value = 0.2532299986575496
$ws.Range("P4").Value = 0.2532299986575496
$ws.Range("Q4").Value = 106.2557772117794
$ws.Range("R4").Value = 956.301994906015
$ws.Range("S4").Value = 0.01227362789561428
$ws.Range("T4").Value = 0.01227362789561428

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 24.50527833333333
$ws.Range("H5").Value = 73.515835
$ws.Range("I5").Value = 0.04846830138877924
$ws.Range("J5").Value = 0.04846830138877924
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 8.941088000000001
$ws.Range("N5").Value = 26.823264
$ws.Range("O5").Value = 0.5221708325715213
$ws.Range("P5").Value = 0.5221708325715213
$ws.Range("Q5").Value = 219.1038500428267
$ws.Range("R5").Value = 1971.93465038544
$ws.Range("S5").Value = 0.02530873328950628
$ws.Range("T5").Value = 0.02530873328950628

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 110.1980973333333
$ws.Range("H6").Value = 330.594292
$ws.Range("I6").Value = 0.2179577200213544
$ws.Range("J6").Value = 0.2179577200213544
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.383699
$ws.Range("N6").Value = 4.151097
$ws.Range("O6").Value = 0.08080976933214185
$ws.Range("P6").Value = 0.08080976933214185
$ws.Range("Q6").Value = 152.480997082036
$ws.Range("R6").Value = 1372.328973738324
$ws.Range("S6").Value = 0.01761311307908521
$ws.Range("T6").Value = 0.01761311307908521

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 110.1980973333333
$ws.Range("H7").Value = 330.594292
$ws.Range("I7").Value = 0.2179577200213544
$ws.Range("J7").Value = 0.2179577200213544
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 2.462094
$ws.Range("N7").Value = 7.386282
$ws.Range("O7").Value = 0.1437893994387872
$ws.Range("P7").Value = 0.1437893994387872
$ws.Range("Q7").Value = 271.318074255816
$ws.Range("R7").Value = 2441.862668302344
$ws.Range("S7").Value = 0.03134000966491788
$ws.Range("T7").Value = 0.03134000966491789

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 110.1980973333333
$ws.Range("H8").Value = 330.594292
$ws.Range("I8").Value = 0.2179577200213544
$ws.Range("J8").Value = 0.2179577200213544
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 4.336036333333333
$ws.Range("N8").Value = 13.008109
$ws.Range("O8").Value = 0.2532299986575496
$ws.Range("P8").Value = 0.2532299986575496
$ws.Range("Q8").Value = 477.8229539015364
$ws.Range("R8").Value = 4300.406585113828
$ws.Range("S8").Value = 0.05519343314841017
$ws.Range("T8").Value = 0.05519343314841017

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 110.1980973333333
$ws.Range("H9").Value = 330.594292
$ws.Range("I9").Value = 0.2179577200213544
$ws.Range("J9").Value = 0.2179577200213544
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 8.941088000000001
$ws.Range("N9").Value = 26.823264
$ws.Range("O9").Value = 0.5221708325715213
$ws.Range("P9").Value = 0.5221708325715213
$ws.Range("Q9").Value = 985.2908856898988
$ws.Range("R9").Value = 8867.617971209089
$ws.Range("S9").Value = 0.1138111641289412
$ws.Range("T9").Value = 0.1138111641289412

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 351.7202226666666
$ws.Range("H10").Value = 1055.160668
$ws.Range("I10").Value = 0.6956575446665283
$ws.Range("J10").Value = 0.6956575446665284
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.383699
$ws.Range("N10").Value = 4.151097
$ws.Range("O10").Value = 0.08080976933214185
$ws.Range("P10").Value = 0.08080976933214185
$ws.Range("Q10").Value = 486.674920383644
$ws.Range("R10").Value = 4380.074283452796
$ws.Range("S10").Value = 0.05621592571866631
$ws.Range("T10").Value = 0.05621592571866632

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 351.7202226666666
$ws.Range("H11").Value = 1055.160668
$ws.Range("I11").Value = 0.6956575446665283
$ws.Range("J11").Value = 0.6956575446665284
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 2.462094
$ws.Range("N11").Value = 7.386282
$ws.Range("O11").Value = 0.1437893994387872
$ws.Range("P11").Value = 0.1437893994387872
$ws.Range("Q11").Value = 865.968249906264
$ws.Range("R11").Value = 7793.714249156376
$ws.Range("S11").Value = 0.1000281805626614
$ws.Range("T11").Value = 0.1000281805626614

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 351.7202226666666
$ws.Range("H12").Value = 1055.160668
$ws.Range("I12").Value = 0.6956575446665283
$ws.Range("J12").Value = 0.6956575446665284
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 4.336036333333333
$ws.Range("N12").Value = 13.008109
$ws.Range("O12").Value = 0.2532299986575496
$ws.Range("P12").Value = 0.2532299986575496
$ws.Range("Q12").Value = 1525.071664650757
$ws.Range("R12").Value = 13725.64498185681
$ws.Range("S12").Value = 0.1761613591020192
$ws.Range("T12").Value = 0.1761613591020193

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 351.7202226666666
$ws.Range("H13").Value = 1055.160668
$ws.Range("I13").Value = 0.6956575446665283
$ws.Range("J13").Value = 0.6956575446665284
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 8.941088000000001
$ws.Range("N13").Value = 26.823264
$ws.Range("O13").Value = 0.5221708325715213
$ws.Range("P13").Value = 0.5221708325715213
$ws.Range("Q13").Value = 3144.761462242261
$ws.Range("R13").Value = 28302.85316018036
$ws.Range("S13").Value = 0.3632520792831814
$ws.Range("T13").Value = 0.3632520792831814

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 19.17031833333333
$ws.Range("H14").Value = 57.510955
$ws.Range("I14").Value = 0.03791643392333802
$ws.Range("J14").Value = 0.03791643392333802
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 1.383699
$ws.Range("N14").Value = 4.151097
$ws.Range("O14").Value = 0.08080976933214185
$ws.Range("P14").Value = 0.08080976933214185
$ws.Range("Q14").Value = 26.525950307515
$ws.Range("R14").Value = 238.733552767635
$ws.Range("S14").Value = 0.003064018279242343
$ws.Range("T14").Value = 0.003064018279242343

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 19.17031833333333
$ws.Range("H15").Value = 57.510955
$ws.Range("I15").Value = 0.03791643392333802
$ws.Range("J15").Value = 0.03791643392333802
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 2.462094
$ws.Range("N15").Value = 7.386282
$ws.Range("O15").Value = 0.1437893994387872
$ws.Range("P15").Value = 0.1437893994387872
$ws.Range("Q15").Value = 47.19912574659
$ws.Range("R15").Value = 424.79213171931
$ws.Range("S15").Value = 0.005451981262697232
$ws.Range("T15").Value = 0.005451981262697233

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 19.17031833333333
$ws.Range("H16").Value = 57.510955
$ws.Range("I16").Value = 0.03791643392333802
$ws.Range("J16").Value = 0.03791643392333802
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 4.336036333333333
$ws.Range("N16").Value = 13.008109
$ws.Range("O16").Value = 0.2532299986575496
$ws.Range("P16").Value = 0.2532299986575496
$ws.Range("Q16").Value = 83.12319681489944
$ws.Range("R16").Value = 748.1087713340951
$ws.Range("S16").Value = 0.009601578511505955
$ws.Range("T16").Value = 0.009601578511505955

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 19.17031833333333
$ws.Range("H17").Value = 57.510955
$ws.Range("I17").Value = 0.03791643392333802
$ws.Range("J17").Value = 0.03791643392333802
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 8.941088000000001
$ws.Range("N17").Value = 26.823264
$ws.Range("O17").Value = 0.5221708325715213
$ws.Range("P17").Value = 0.5221708325715213
$ws.Range("Q17").Value = 171.4035032063467
$ws.Range("R17").Value = 1542.63152885712
$ws.Range("S17").Value = 0.01979885586989249
$ws.Range("T17").Value = 0.01979885586989249
